$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 21.285715
$ws.Range("I11").Value = 21.285715
$ws.Range("K11").Value = 21.285715
$ws.Range("M11").Value = 118.714285

$ws.Range("H132").Value = 6537271.5
$ws.Range("I132").Value = 6803833.5
$ws.Range("K132").Value = 20411500.5
$ws.Range("M132").Value = -20408970.5

$ws.Range("H137").Value = 2318.1333
$ws.Range("I137").Value = 2177.1
$ws.Range("K137").Value = 6531.299999999999
$ws.Range("M137").Value = -3981.299999999999

$ws.Range("H138").Value = 3836.5095
$ws.Range("I138").Value = 3041.25
$ws.Range("J138").Value = 3977.889
$ws.Range("K138").Value = 9123.75
$ws.Range("L138").Value = 11933.667
$ws.Range("M138").Value = -3983.75
$ws.Range("N138").Value = -22213.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2781.25
$ws.Range("I132").Value = 2189.1
$ws.Range("J132").Value = 5742
$ws.Range("K132").Value = 6567.299999999999
$ws.Range("L132").Value = 17226
$ws.Range("M132").Value = -4037.299999999999
$ws.Range("N132").Value = -22286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1785.8462
$ws.Range("I105").Value = 1571.6
$ws.Range("K105").Value = 1571.6
$ws.Range("M105").Value = 175.4000000000001

$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H133").Value = 71743
$ws.Range("J133").Value = 85657.336
$ws.Range("L133").Value = 85657.336
$ws.Range("N133").Value = -95777.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 208.4
$ws.Range("J6").Value = 210.5
$ws.Range("L6").Value = 210.5
$ws.Range("N6").Value = -436.5

$ws.Range("H15").Value = 2557.2222
$ws.Range("J15").Value = 2839.625
$ws.Range("L15").Value = 2839.625
$ws.Range("N15").Value = -3179.625

$ws.Range("H16").Value = 12500
$ws.Range("I16").Value = 5000
$ws.Range("K16").Value = 5000
$ws.Range("M16").Value = -4713

$ws.Range("H31").Value = 52204.57
$ws.Range("I31").Value = 78727.69500000001
$ws.Range("K31").Value = 78727.69500000001
$ws.Range("M31").Value = -78432.69500000001

$ws.Range("H34").Value = 52204.57
$ws.Range("I34").Value = 78727.69500000001
$ws.Range("K34").Value = 78727.69500000001
$ws.Range("M34").Value = -78525.69500000001

$ws.Range("H58").Value = 3341.524
$ws.Range("I58").Value = 3851.0667
$ws.Range("J58").Value = 2067.6667
$ws.Range("K58").Value = 3851.0667
$ws.Range("L58").Value = 2067.6667
$ws.Range("M58").Value = -3648.0667
$ws.Range("N58").Value = -2473.6667

$ws.Range("H62").Value = 4710.9
$ws.Range("I62").Value = 3800
$ws.Range("J62").Value = 5101.2856
$ws.Range("K62").Value = 3800
$ws.Range("L62").Value = 5101.2856
$ws.Range("M62").Value = -3176
$ws.Range("N62").Value = -6349.2856

$ws.Range("H65").Value = 4710.9
$ws.Range("I65").Value = 3800
$ws.Range("J65").Value = 5101.2856
$ws.Range("K65").Value = 19000
$ws.Range("L65").Value = 25506.428
$ws.Range("M65").Value = -15880
$ws.Range("N65").Value = -31746.428

$ws.Range("H105").Value = 1601.4
$ws.Range("I105").Value = 1501.75
$ws.Range("K105").Value = 1501.75
$ws.Range("M105").Value = 245.25

$ws.Range("H107").Value = 717.25
$ws.Range("I107").Value = 324.5
$ws.Range("K107").Value = 324.5
$ws.Range("M107").Value = 1595.5

$ws.Range("H113").Value = 12500
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

$ws.Range("H136").Value = 3341.524
$ws.Range("I136").Value = 3851.0667
$ws.Range("J136").Value = 2067.6667
$ws.Range("K136").Value = 11553.2001
$ws.Range("L136").Value = 6203.000100000001
$ws.Range("M136").Value = -9003.2001
$ws.Range("N136").Value = -11303.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2200
$ws.Range("J39").Value = 2200
$ws.Range("L39").Value = 6600
$ws.Range("N39").Value = -7188

$ws.Range("H40").Value = 335.6316
$ws.Range("I40").Value = 196.6923
$ws.Range("K40").Value = 786.7692
$ws.Range("M40").Value = -717.7692

$ws.Range("H41").Value = 724.9091
$ws.Range("I41").Value = 40
$ws.Range("J41").Value = 877.1111
$ws.Range("K41").Value = 120
$ws.Range("L41").Value = 2631.3333
$ws.Range("M41").Value = 218
$ws.Range("N41").Value = -3307.3333

$ws.Range("H86").Value = 372
$ws.Range("J86").Value = 372
$ws.Range("L86").Value = 1116
$ws.Range("N86").Value = -3488

$ws.Range("H89").Value = 372
$ws.Range("J89").Value = 372
$ws.Range("L89").Value = 3348
$ws.Range("N89").Value = -15204

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 266.15625
$ws.Range("I2").Value = 255.2381
$ws.Range("K2").Value = 255.2381
$ws.Range("M2").Value = -142.2381

$ws.Range("H80").Value = 5054.625
$ws.Range("I80").Value = 3724.9092
$ws.Range("K80").Value = 3724.9092
$ws.Range("M80").Value = -2726.9092

$ws.Range("H83").Value = 5054.625
$ws.Range("I83").Value = 3724.9092
$ws.Range("K83").Value = 18624.546
$ws.Range("M83").Value = -13632.546

$ws.Range("H97").Value = 1121.4117
$ws.Range("I97").Value = 836.0417
$ws.Range("J97").Value = 1806.3
$ws.Range("K97").Value = 836.0417
$ws.Range("L97").Value = 1806.3
$ws.Range("M97").Value = -340.0417
$ws.Range("N97").Value = -2798.3

$ws.Range("H122").Value = 2085.6365
$ws.Range("I122").Value = 1643.6666
$ws.Range("K122").Value = 4930.9998
$ws.Range("M122").Value = -2480.9998

$ws.Range("H132").Value = 3844.6667
$ws.Range("I132").Value = 3608.5715
$ws.Range("K132").Value = 10825.7145
$ws.Range("M132").Value = -8295.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9299.723
$ws.Range("I7").Value = 10536.417
$ws.Range("K7").Value = 10536.417
$ws.Range("M7").Value = -10424.417

$ws.Range("H16").Value = 38467844
$ws.Range("I16").Value = 71432100
$ws.Range("J16").Value = 9549.833000000001
$ws.Range("K16").Value = 71432100
$ws.Range("L16").Value = 9549.833000000001
$ws.Range("M16").Value = -71431930
$ws.Range("N16").Value = -9889.833000000001

$ws.Range("H40").Value = 4050.1155
$ws.Range("I40").Value = 2990.25
$ws.Range("K40").Value = 2990.25
$ws.Range("M40").Value = -2854.25

$ws.Range("H46").Value = 1952.625
$ws.Range("I46").Value = 1243
$ws.Range("J46").Value = 3371.875
$ws.Range("K46").Value = 1243
$ws.Range("L46").Value = 3371.875
$ws.Range("M46").Value = -1055
$ws.Range("N46").Value = -3747.875

$ws.Range("H93").Value = 2758
$ws.Range("I93").Value = 3965
$ws.Range("K93").Value = 3965
$ws.Range("M93").Value = -2717

$ws.Range("H126").Value = 9299.723
$ws.Range("I126").Value = 10536.417
$ws.Range("K126").Value = 31609.251
$ws.Range("M126").Value = -29139.251

$ws.Range("H132").Value = 4982.375
$ws.Range("I132").Value = 4474.8184
$ws.Range("K132").Value = 13424.4552
$ws.Range("M132").Value = -10894.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H100").Value = 1427.7333
$ws.Range("I100").Value = 838.7143
$ws.Range("K100").Value = 1677.4286
$ws.Range("M100").Value = -1136.4286

$ws.Range("H113").Value = 956.2941
$ws.Range("I113").Value = 528.6667
$ws.Range("K113").Value = 1586.0001
$ws.Range("M113").Value = 583.9999

$ws.Range("H132").Value = 5413.484
$ws.Range("I132").Value = 4958.222
$ws.Range("K132").Value = 14874.666
$ws.Range("M132").Value = -12344.666

$ws.Range("H136").Value = 3135.9092
$ws.Range("I136").Value = 2969
$ws.Range("K136").Value = 8907
$ws.Range("M136").Value = -6357
